# Apply the edits described by the commit:
# "Add python script to implement progress bar for timeshift on scale of
#  OSD slider. Tidy up timeshift display in seek bar."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status/comment for the OSD timeshift progress row (row 41):
# Status moves on from "Not skinnable" to "Partially done", and the comment
# now describes the python script work-around plus the outstanding bug.
$ws.Range("B41").Value = "Partially done"
$ws.Range("C41").Value = "Done with Python script. Still awaiting bug fix for Player.Progress (bug tracker #17469)"

# Fix typo in the "Not skinnable" explanation for the context-menu row (row 40):
# "ca be added" -> "can be added"
$ws.Range("C40").Value = "There is very limited scope for skins to do this (existing context menus can be added to, but new ones can't be created)"

# Update the view/selection to reflect where the author ended up editing.
$ws.Range("C40").Select()
